# Auto-generated Excel COM-interop edit script
# Applies updated crypto price/volume figures (and one coin swap at row 51)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextCell($cellRef, $text) {
    $rng = $ws.Range($cellRef)
    $rng.NumberFormat = '@'
    $rng.Value = $text
    $rng.Style = 'Normal'
}

Set-TextCell 'D2' '28.911.42'
Set-TextCell 'E2' '  -1.66%  '
Set-TextCell 'D3' '1.833.86'
Set-TextCell 'E3' '  -1.90%  '
Set-TextCell 'D4' '0.9994'
Set-TextCell 'E4' '  -0.15%  '
Set-TextCell 'D5' '244.60'
Set-TextCell 'E5' '  +0.37%  '
Set-TextCell 'D6' '0.6896'
Set-TextCell 'E6' '  -2.03%  '
Set-TextCell 'D8' '0.07688'
Set-TextCell 'E8' '  -2.99%  '
Set-TextCell 'D9' '0.3050'
Set-TextCell 'E9' '  -2.62%  '
Set-TextCell 'D10' '23.40'
Set-TextCell 'E10' '  -4.43%  '
Set-TextCell 'D11' '0.07779'
Set-TextCell 'E11' '  -0.76%  '
Set-TextCell 'D12' '1.826.94'
Set-TextCell 'E12' '  -3.92%  '
Set-TextCell 'D13' '5.077'
Set-TextCell 'E13' '  -1.80%  '
Set-TextCell 'D14' '90.44'
Set-TextCell 'E14' '  -3.54%  '
Set-TextCell 'D15' '0.6806'
Set-TextCell 'E15' '  -2.87%  '
Set-TextCell 'D16' '6.445'
Set-TextCell 'E16' '  -1.13%  '
Set-TextCell 'D17' '0.000008288'
Set-TextCell 'E17' '  -1.39%  '
Set-TextCell 'D18' '28.907.35'
Set-TextCell 'E18' '  -1.96%  '
Set-TextCell 'D19' '242.83'
Set-TextCell 'E19' '  -3.78%  '
Set-TextCell 'D20' '2.077.52'
Set-TextCell 'E20' '  -3.68%  '
Set-TextCell 'E21' '  -2.99%  '
Set-TextCell 'D22' '0.9996'
Set-TextCell 'E22' '  -0.09%  '
Set-TextCell 'D23' '7.476'
Set-TextCell 'E23' '  -2.50%  '
Set-TextCell 'D24' '1.000'
Set-TextCell 'E24' '  -0.09%  '
Set-TextCell 'D25' '162.47'
Set-TextCell 'E25' '  +0.48%  '
Set-TextCell 'D26' '0.1472'
Set-TextCell 'E26' '  -5.19%  '
Set-TextCell 'D27' '8.808'
Set-TextCell 'E27' '  -2.31%  '
Set-TextCell 'D28' '18.19'
Set-TextCell 'D29' '1.546'
Set-TextCell 'E29' '  +2.63%  '
Set-TextCell 'D30' '4.209'
Set-TextCell 'E30' '  -2.49%  '
Set-TextCell 'D31' '4.148'
Set-TextCell 'E31' '  -2.58%  '
Set-TextCell 'D32' '1.184'
Set-TextCell 'E32' '  -2.54%  '
Set-TextCell 'D33' '0.05105'
Set-TextCell 'E33' '  -3.08%  '
Set-TextCell 'D34' '0.7639'
Set-TextCell 'E34' '  +1.67%  '
Set-TextCell 'E35' '  -2.75%  '
Set-TextCell 'E36' '  -3.27%  '
Set-TextCell 'D37' '2.686'
Set-TextCell 'E37' '  -0.89%  '
Set-TextCell 'E38' '  -1.54%  '
Set-TextCell 'D39' '1.221.70'
Set-TextCell 'E39' '  -3.94%  '
Set-TextCell 'D40' '2.697'
Set-TextCell 'E40' '  -2.62%  '
Set-TextCell 'D41' '0.9399'
Set-TextCell 'E41' '  +5.32%  '
Set-TextCell 'D42' '108.20'
Set-TextCell 'E42' '  -1.01%  '
Set-TextCell 'E43' '  -0.11%  '
Set-TextCell 'D44' '5.677'
Set-TextCell 'E44' '  -5.97%  '
Set-TextCell 'D45' '9.661'
Set-TextCell 'E45' '  +0.31%  '
Set-TextCell 'D46' '0.00000000122'
Set-TextCell 'E46' '  -2.65%  '
Set-TextCell 'D47' '0.5172'
Set-TextCell 'E47' '  -0.25%  '
Set-TextCell 'D48' '1.977.04'
Set-TextCell 'E48' '  -3.26%  '
Set-TextCell 'D49' '64.14'
Set-TextCell 'E49' '  -9.51%  '
Set-TextCell 'D50' '1.747'
Set-TextCell 'E50' '  -3.14%  '
Set-TextCell 'B51' 'TheSandbox'
Set-TextCell 'C51' 'https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand'
Set-TextCell 'D51' '0.4181'
Set-TextCell 'E51' '  -2.79%  '
